# corrected index problem with excel
# Columns E, J, O ("Free Lies", "Free Lies Sec.", "# Free Lies") were
# mis-aligned by two rows relative to the rest of the table (columns
# A-D, G-I, L-N). This re-aligns the label column (A) for rows 8-13 and
# shifts the E/J/O values down into their correct rows (8-17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Participants) relabeling, rows 8-13 ---
$ws.Range("A8").Value  = "F27_19"
$ws.Range("A9").Value  = "F29_14"
$ws.Range("A10").Value = "M24_18"
$ws.Range("A11").Value = "M25_5"
$ws.Range("A12").Value = "M26_10"
$ws.Range("A13").Value = "M26_12"

# --- Column E (Free Lies) ---
$ws.Range("E8").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("E10").Value = 0.1469600724073327
$ws.Range("E11").ClearContents()
$ws.Range("E12").Value = 0.3216824994606119
$ws.Range("E13").Value = -0.1409664951635983
$ws.Range("E14").Value = 0.1928504054064898
$ws.Range("E15").Value = -0.03326235607924657
$ws.Range("E17").Value = 0.08146840787291081

# --- Column J (Free Lies Sec.) ---
$ws.Range("J8").ClearContents()
$ws.Range("J9").ClearContents()
$ws.Range("J10").Value = 3.61353568
$ws.Range("J11").ClearContents()
$ws.Range("J12").Value = 3.270023816666667
$ws.Range("J13").Value = 4.22070955
$ws.Range("J14").Value = 3.717797666666666
$ws.Range("J15").Value = 7.309348783333333
$ws.Range("J17").Value = 3.6195893

# --- Column O (# Free Lies) ---
$ws.Range("O8").ClearContents()
$ws.Range("O9").ClearContents()
$ws.Range("O10").Value = 5
$ws.Range("O11").ClearContents()
$ws.Range("O14").Value = 6
$ws.Range("O15").Value = 6
$ws.Range("O17").Value = 2
